# Fixed the emissions from RD distribution for algae CAP and HTL pathways,
# also changed the RD distribution loss factor on the "Fuel dist urban" sheet.
# Column C (Renewable Diesel) is the entered input; columns D:G are formulas
# (=C<row>) that ripple automatically once recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fuel dist urban")

$newValues = @{
    4  = [double]"5093.7824552240209"
    5  = [double]"5023.1958652604862"
    6  = [double]"41.657173986229949"
    7  = [double]"527.78573845908159"
    8  = [double]"4453.7529528151745"
    9  = [double]"9.869537321209304E-2"
    10 = [double]"1.1883133447180589"
    11 = [double]"0.60557378989801824"
    12 = [double]"1.3472060329660924"
    13 = [double]"4.077408879999693E-2"
    14 = [double]"2.866035992117464E-2"
    15 = [double]"2.4510023670013949E-2"
    16 = [double]"3.3183046022218594E-3"
    17 = [double]"1.5425727512154881E-2"
    18 = [double]"0.48940655600260441"
    19 = [double]"5.8402637809921217E-3"
    20 = [double]"387.85040331686662"
    21 = [double]"0.26551552589021865"
    22 = [double]"0.17288665151109295"
    23 = [double]"0.21330355735173617"
    24 = [double]"9.4740252492059081E-3"
    25 = [double]"5.0442745092267195E-3"
    26 = [double]"5.2423279836688944E-3"
    27 = [double]"5.4229333446464761E-4"
    28 = [double]"1.8571109385136862E-3"
}

foreach ($row in $newValues.Keys) {
    $ws.Range("C$row").Value = $newValues[$row]
}

# Force recalculation so the dependent D:G formula cells (=C<row>) pick up
# the new cached values (workbook calc mode is manual).
$excel.Calculate()

# Restore the view state: "Fuel dist urban" becomes the active/selected
# sheet with C10 selected; "Fuel specs" is no longer the selected tab.
$ws.Activate()
$ws.Range("C10").Select()

Write-Host "Updated Renewable Diesel (column C) values on 'Fuel dist urban' for rows 4-28."
